$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D3').Value = 44426
$ws.Range('H3').Value = 'Española'
$ws.Range('J3').Value = 600
$ws.Range('K3').Value = 11500
$ws.Range('L3').Value = 12000
$ws.Range('M3').Value = 11750
$ws.Range('N3').Value = '$/caja 30 unidades'
$ws.Range('O3').Value = 'Provincia de Limarí'
$ws.Range('P3').Value = 392
$ws.Range('Q3').Value = 30

$ws.Range('D4').Value = 44426
$ws.Range('H4').Value = 'Madrigal'
$ws.Range('J4').Value = 500
$ws.Range('K4').Value = 12500
$ws.Range('L4').Value = 13000
$ws.Range('M4').Value = 12750
$ws.Range('N4').Value = '$/caja 40 unidades'
$ws.Range('O4').Value = 'Provincia de Limarí'
$ws.Range('P4').Value = 319
$ws.Range('Q4').Value = 40

$ws.Range('D5').Value = 44701
$ws.Range('H5').Value = 'Española'
$ws.Range('J5').Value = 400
$ws.Range('K5').Value = 19000
$ws.Range('L5').Value = 20000
$ws.Range('M5').Value = 19500
$ws.Range('N5').Value = '$/caja 30 unidades'
$ws.Range('O5').Value = 'Provincia del Elquí'
$ws.Range('P5').Value = 650
$ws.Range('Q5').Value = 30

$ws.Range('D6').Value = 44498
$ws.Range('H6').Value = 'Española'
$ws.Range('J6').Value = 400
$ws.Range('K6').Value = 8500
$ws.Range('L6').Value = 9000
$ws.Range('M6').Value = 8750
$ws.Range('N6').Value = '$/caja 30 unidades'
$ws.Range('O6').Value = 'Provincia de Limarí'
$ws.Range('P6').Value = 292
$ws.Range('Q6').Value = 30

$ws.Range('D7').Value = 44687
$ws.Range('H7').Value = 'Española'
$ws.Range('J7').Value = 400
$ws.Range('K7').Value = 18000
$ws.Range('L7').Value = 19000
$ws.Range('M7').Value = 18500
$ws.Range('N7').Value = '$/caja 30 unidades'
$ws.Range('O7').Value = 'Provincia de Limarí'
$ws.Range('P7').Value = 617
$ws.Range('Q7').Value = 30

$ws.Range('D8').Value = 44420
$ws.Range('H8').Value = 'Madrigal'
$ws.Range('J8').Value = 800
$ws.Range('K8').Value = 14000
$ws.Range('L8').Value = 15000
$ws.Range('M8').Value = 14500
$ws.Range('N8').Value = '$/caja 40 unidades'
$ws.Range('O8').Value = 'Provincia de Limarí'
$ws.Range('P8').Value = 362
$ws.Range('Q8').Value = 40

$ws.Range('D9').Value = 44420
$ws.Range('H9').Value = 'Madrigal'
$ws.Range('J9').Value = 700
$ws.Range('K9').Value = 13000
$ws.Range('L9').Value = 14000
$ws.Range('M9').Value = 13500
$ws.Range('N9').Value = '$/caja 40 unidades'
$ws.Range('O9').Value = 'Provincia del Elquí'
$ws.Range('P9').Value = 338
$ws.Range('Q9').Value = 40

$ws.Range('D10').Value = 44427
$ws.Range('H10').Value = 'Madrigal'
$ws.Range('J10').Value = 400
$ws.Range('K10').Value = 12000
$ws.Range('L10').Value = 13000
$ws.Range('M10').Value = 12500
$ws.Range('N10').Value = '$/caja 40 unidades'
$ws.Range('O10').Value = 'Provincia de Limarí'
$ws.Range('P10').Value = 312
$ws.Range('Q10').Value = 40

$ws.Range('D11').Value = 44484
$ws.Range('H11').Value = 'Española'
$ws.Range('J11').Value = 300
$ws.Range('K11').Value = 9000
$ws.Range('L11').Value = 10000
$ws.Range('M11').Value = 9500
$ws.Range('N11').Value = '$/caja 30 unidades'
$ws.Range('O11').Value = 'Provincia del Elquí'
$ws.Range('P11').Value = 317
$ws.Range('Q11').Value = 30

$ws.Range('D12').Value = 44438
$ws.Range('H12').Value = 'Española'
$ws.Range('J12').Value = 400
$ws.Range('K12').Value = 11000
$ws.Range('L12').Value = 12000
$ws.Range('M12').Value = 11500
$ws.Range('N12').Value = '$/caja 30 unidades'
$ws.Range('O12').Value = 'Provincia del Elquí'
$ws.Range('P12').Value = 383
$ws.Range('Q12').Value = 30
